# Fixed bug in Model_patients using wrong data to test model
# Adds a "MEAN" summary column to both the Precision table (left) and the
# Recall table (right), shifting the Recall table one column to the right
# to keep a blank spacer column between the two tables.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the second (Recall) table one column to the right (H:M -> I:N),
# leaving column H as a new blank spacer and freeing up column G for the
# first table's new MEAN column.
$ws.Columns("H").Insert()

# Give the new blank spacer cell G1 the same look as the rest of row 1
# (bold, centered) by copying the formatting from the cell beside it.
$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").HorizontalAlignment = $ws.Range("F1").HorizontalAlignment

# --- MEAN column for the first (Precision) table --------------------------
$ws.Range("G2").Value = "MEAN"
$ws.Range("G2").Font.Bold = $true
$ws.Range("G2").Font.Italic = $true

$ws.Range("G3").Formula = "=AVERAGE(B3:F3)"
$ws.Range("G4").Formula = "=AVERAGE(B4:F4)"
$ws.Range("G5").Formula = "=AVERAGE(B5:F5)"
$ws.Range("G6").Formula = "=AVERAGE(B6:F6)"

$ws.Range("G3:G6").NumberFormat = "0.00"
$ws.Range("G3:G5").Font.Italic = $true
$ws.Range("G6").Font.Bold = $true
$ws.Range("G6").Font.Italic = $true

# --- MEAN column for the second (Recall) table -----------------------------
$ws.Range("O2").Value = "MEAN"
$ws.Range("O2").Font.Bold = $true
$ws.Range("O2").Font.Italic = $true

$ws.Range("O3").Formula = "=AVERAGE(J3:N3)"
$ws.Range("O4").Formula = "=AVERAGE(J4:N4)"
$ws.Range("O5").Formula = "=AVERAGE(J5:N5)"
$ws.Range("O6").Formula = "=AVERAGE(J6:N6)"

$ws.Range("O3:O6").NumberFormat = "0.00"
$ws.Range("O3:O5").Font.Italic = $true
$ws.Range("O6").Font.Bold = $true
$ws.Range("O6").Font.Italic = $true

# --- Extend / recreate the data-bar conditional formatting -----------------
# Widen the Precision table's data bar to cover the new MEAN column.
$fcPrecision = $ws.Range("B3:F6").FormatConditions.Item(1)
$fcPrecision.ModifyAppliesToRange($ws.Range("B3:G6"))

# Move the Recall table's data bar to its new location (shifted right).
$fcRecall = $ws.Range("I3:M6").FormatConditions.Item(1)
$fcRecall.ModifyAppliesToRange($ws.Range("J3:N6"))

# Add a fresh data bar for the new Recall-table MEAN column and give it the
# highest priority, matching the original rule ordering.
$fcMean = $ws.Range("O3:O6").FormatConditions.AddDatabar()
$fcMean.SetFirstPriority()

# --- Restore the active selection ------------------------------------------
$ws.Range("F5").Select()
